# Fruta / hortaliza, semanal
# Update the weekly price records for "Vega Monumental Concepción - Tuna".
# Only the Fecha (D), Calidad (L), Volumen (M), Precio mínimo (N),
# Precio máximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# columns change; everything else on each row stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44280
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806

# Row 3
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 667

# Row 4
$ws.Range("D4").Value = 44699
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21000
$ws.Range("S4").Value = 1167

# Row 5
$ws.Range("D5").Value = 44699
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44516
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 33000
$ws.Range("O6").Value = 34000
$ws.Range("P6").Value = 33500
$ws.Range("S6").Value = 1861

# Row 7
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 1111

# Row 8
$ws.Range("D8").Value = 44687
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("S8").Value = 1028
